# Insert a new "orphaCase_no" column right after "orphaCoding_no" (column C),
# shifting unique_rdCase_no / rdCase_no / case_no / patient_no / inpatientCases_no
# one column to the right (D->E, E->F, F->G, G->H, H->I), and populate the new
# column's header + data value for the existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at D; everything from D onward shifts right to E..I.
$ws.Columns.Item(4).Insert()

# Fill in the newly inserted column D.
$ws.Range("D1").Value = "orphaCase_no"
$ws.Range("D2").Value = 124
